$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '44.222.27'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '2.233.75'
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('E4').Value = '  +0.61%  '
$ws.Range('D5').Value = '306.25'
$ws.Range('E5').Value = '  -3.44%  '
$ws.Range('D6').Value = '93.21'
$ws.Range('E6').Value = '  -6.94%  '
$ws.Range('E7').Value = '  -1.39%  '
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  -3.73%  '
$ws.Range('D10').Value = '33.94'
$ws.Range('E10').Value = '  -6.53%  '
$ws.Range('D11').Value = '0.0805'
$ws.Range('E11').Value = '  -2.78%  '
$ws.Range('D12').Value = '7.09'
$ws.Range('E12').Value = '  -4.65%  '
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '2.320.60'
$ws.Range('E14').Value = '  +3.48%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '0.822'
$ws.Range('E15').Value = '  -3.47%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '13.36'
$ws.Range('E16').Value = '  -5.07%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '44.102.63'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.0₃0959'
$ws.Range('E18').Value = '  -2.79%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '6.31'
$ws.Range('E19').Value = '  -1.24%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').Value = '11.92'
$ws.Range('E20').Value = '  -10.22%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').Value = '65.43'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('B22').Value = 'PancakeSwap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D22').Value = '3.08'
$ws.Range('E22').Value = '  +2.54%  '
$ws.Range('D23').Value = '236.31'
$ws.Range('E23').Value = '  -1.76%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').Value = '1.98'
$ws.Range('E24').Value = '  -3.31%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('B26').Value = 'InjectiveProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D26').Value = '40.02'
$ws.Range('E26').Value = '  +4.56%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').Value = '2.20'
$ws.Range('E27').Value = '  +3.68%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '9.79'
$ws.Range('E28').Value = '  -3.99%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '19.96'
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '5.79'
$ws.Range('E30').Value = '  -4.54%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = '151.79'
$ws.Range('E31').Value = '  -2.87%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.0788'
$ws.Range('E32').Value = '  -6.64%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').Value = '2.60'
$ws.Range('E33').Value = '  -2.59%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '3.04'
$ws.Range('E34').Value = '  -12.34%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').Value = '0.119'
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '0.108'
$ws.Range('E36').Value = '  -4.04%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '1.73'
$ws.Range('E37').Value = '  -10.43%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = '3.44'
$ws.Range('E38').Value = '  -3.28%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').Value = '14.24'
$ws.Range('E39').Value = '  -7.58%  '
$ws.Range('D40').Value = '3.76'
$ws.Range('E40').Value = '  -4.53%  '
$ws.Range('D41').Value = '0.0296'
$ws.Range('E41').Value = '  -4.56%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '1.01'
$ws.Range('E42').Value = '  +0.44%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.706.82'
$ws.Range('E43').Value = '  -1.27%  '
$ws.Range('B44').Value = 'BitcoinSV'
$ws.Range('C44').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D44').Value = '81.98'
$ws.Range('E44').Value = '  -2.88%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = '0.190'
$ws.Range('E45').Value = '  -3.94%  '
$ws.Range('B46').Value = 'THORChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D46').Value = '4.90'
$ws.Range('E46').Value = '  -6.22%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '98.54'
$ws.Range('E47').Value = '  -4.09%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = '1.59'
$ws.Range('E48').Value = '  -2.66%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = '54.17'
$ws.Range('E49').Value = '  -5.02%  '
$ws.Range('D50').Value = '7.98'
$ws.Range('E50').Value = '  -2.34%  '
$ws.Range('B51').Value = 'ordi'
$ws.Range('C51').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D51').Value = '66.12'
$ws.Range('E51').Value = '  -6.74%  '
